# Auto-generated edit script applying the Cerberus_Profits.xlsx diff
# (market-price / profit recalculation across the 8 crafting-job sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # ALC
$ws.Range("H28").Value = 2722.4
$ws.Range("I28").Value = 4520.8335
$ws.Range("J28").Value = 1523.4445
$ws.Range("K28").Value = 4520.8335
$ws.Range("L28").Value = 1523.4445
$ws.Range("M28").Value = -4035.8335
$ws.Range("N28").Value = -2493.4445

$ws.Range("H62").Value = 50003160
$ws.Range("I62").Value = 62502976
$ws.Range("K62").Value = 62502976
$ws.Range("M62").Value = -62502352

$ws.Range("H65").Value = 50003160
$ws.Range("I65").Value = 62502976
$ws.Range("K65").Value = 312514880
$ws.Range("M65").Value = -312511760

$ws.Range("H70").Value = 10267.923
$ws.Range("J70").Value = 15123.5
$ws.Range("L70").Value = 45370.5
$ws.Range("N70").Value = -45910.5

$ws.Range("H73").Value = 10267.923
$ws.Range("J73").Value = 15123.5
$ws.Range("L73").Value = 45370.5
$ws.Range("N73").Value = -47242.5

$ws.Range("H74").Value = 6649.92
$ws.Range("I74").Value = 5677.091
$ws.Range("J74").Value = 7414.2856
$ws.Range("K74").Value = 5677.091
$ws.Range("L74").Value = 7414.2856
$ws.Range("M74").Value = -4741.091
$ws.Range("N74").Value = -9286.285599999999

$ws.Range("H77").Value = 6649.92
$ws.Range("I77").Value = 5677.091
$ws.Range("J77").Value = 7414.2856
$ws.Range("K77").Value = 28385.455
$ws.Range("L77").Value = 37071.428
$ws.Range("M77").Value = -23705.455
$ws.Range("N77").Value = -46431.428

$ws.Range("H86").Value = 1856.8572
$ws.Range("I86").Value = 1916.3334
$ws.Range("K86").Value = 1916.3334
$ws.Range("M86").Value = -793.3334

$ws.Range("H89").Value = 1856.8572
$ws.Range("I89").Value = 1916.3334
$ws.Range("K89").Value = 9581.666999999999
$ws.Range("M89").Value = -3965.666999999999

$ws.Range("H107").Value = 1551.5
$ws.Range("I107").Value = 1489.375
$ws.Range("K107").Value = 1489.375
$ws.Range("M107").Value = 430.625

$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("M111").ClearContents()

$ws.Range("H113").Value = 5614.4546
$ws.Range("I113").Value = 4312.2856
$ws.Range("K113").Value = 4312.2856
$ws.Range("M113").Value = -1058.2856

$ws.Range("H131").Value = 4851.1
$ws.Range("I131").Value = 2387.2856
$ws.Range("J131").Value = 10600
$ws.Range("K131").Value = 7161.8568
$ws.Range("L131").Value = 31800
$ws.Range("M131").Value = -2121.8568
$ws.Range("N131").Value = -41880

$ws = $wb.Worksheets.Item(2)  # ARM
$ws.Range("H45").Value = 1969.3334
$ws.Range("I45").Value = 1495.6
$ws.Range("J45").Value = 2561.5
$ws.Range("K45").Value = 1495.6
$ws.Range("L45").Value = 2561.5
$ws.Range("M45").Value = -1118.6
$ws.Range("N45").Value = -3315.5

$ws.Range("H63").Value = 1552.25
$ws.Range("I63").Value = 1552.25
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 1552.25
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -866.25
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 1552.25
$ws.Range("I66").Value = 1552.25
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 7761.25
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -4329.25
$ws.Range("N66").ClearContents()

$ws.Range("H74").Value = 2169.9524
$ws.Range("I74").Value = 1654.8334
$ws.Range("J74").Value = 2856.7778
$ws.Range("K74").Value = 1654.8334
$ws.Range("L74").Value = 2856.7778
$ws.Range("M74").Value = -780.8334
$ws.Range("N74").Value = -4604.7778

$ws.Range("H77").Value = 2169.9524
$ws.Range("I77").Value = 1654.8334
$ws.Range("J77").Value = 2856.7778
$ws.Range("K77").Value = 8274.166999999999
$ws.Range("L77").Value = 14283.889
$ws.Range("M77").Value = -3906.166999999999
$ws.Range("N77").Value = -23019.889

$ws = $wb.Worksheets.Item(3)  # BSM
$ws.Range("H80").Value = 879.4167
$ws.Range("J80").Value = 914.6667
$ws.Range("L80").Value = 914.6667
$ws.Range("N80").Value = -2910.6667

$ws.Range("H83").Value = 879.4167
$ws.Range("J83").Value = 914.6667
$ws.Range("L83").Value = 4573.3335
$ws.Range("N83").Value = -14557.3335

$ws = $wb.Worksheets.Item(4)  # CRP
$ws.Range("H31").Value = 2580.366
$ws.Range("I31").Value = 1648.3334
$ws.Range("J31").Value = 3559
$ws.Range("K31").Value = 1648.3334
$ws.Range("L31").Value = 3559
$ws.Range("M31").Value = -1353.3334
$ws.Range("N31").Value = -4149

$ws.Range("H34").Value = 2580.366
$ws.Range("I34").Value = 1648.3334
$ws.Range("J34").Value = 3559
$ws.Range("K34").Value = 1648.3334
$ws.Range("L34").Value = 3559
$ws.Range("M34").Value = -1446.3334
$ws.Range("N34").Value = -3963

$ws.Range("H99").Value = 2166.5386

$ws.Range("H126").Value = 2166.5386

$ws = $wb.Worksheets.Item(5)  # CUL
$ws.Range("H2").Value = 171.88889
$ws.Range("I2").Value = 360.25
$ws.Range("J2").Value = 21.2
$ws.Range("K2").Value = 2161.5
$ws.Range("L2").Value = 127.2
$ws.Range("M2").Value = -2048.5
$ws.Range("N2").Value = -353.2

$ws.Range("H95").Value = 4999
$ws.Range("J95").Value = 4999
$ws.Range("L95").Value = 14997
$ws.Range("N95").Value = -19115

$ws.Range("H131").Value = 4493040.5
$ws.Range("I131").Value = 6946077
$ws.Range("J131").Value = 4020166.2
$ws.Range("K131").Value = 20838231
$ws.Range("L131").Value = 12060498.6
$ws.Range("M131").Value = -20833191
$ws.Range("N131").Value = -12070578.6

$ws = $wb.Worksheets.Item(6)  # GSM
$ws.Range("H102").Value = 3818.44
$ws.Range("I102").Value = 3842.5
$ws.Range("K102").Value = 3842.5
$ws.Range("M102").Value = -2220.5

$ws.Range("H126").Value = 8626.5
$ws.Range("J126").Value = 8500
$ws.Range("L126").Value = 25500
$ws.Range("N126").Value = -30440

$ws.Range("H141").Value = 89237.836
$ws.Range("J141").Value = 89237.836
$ws.Range("L141").Value = 89237.836
$ws.Range("N141").Value = -99597.836

$ws = $wb.Worksheets.Item(7)  # LTW
$ws.Range("H40").Value = 2346.6155

$ws.Range("H61").Value = 2028.4286
$ws.Range("I61").Value = 2033.1666
$ws.Range("K61").Value = 2033.1666
$ws.Range("M61").Value = -1831.1666

$ws.Range("H68").Value = 2645.125
$ws.Range("I68").Value = 2485.1667
$ws.Range("J68").Value = 3125
$ws.Range("K68").Value = 2485.1667
$ws.Range("L68").Value = 3125
$ws.Range("M68").Value = -1736.1667
$ws.Range("N68").Value = -4623

$ws.Range("H71").Value = 2645.125
$ws.Range("I71").Value = 2485.1667
$ws.Range("J71").Value = 3125
$ws.Range("K71").Value = 12425.8335
$ws.Range("L71").Value = 15625
$ws.Range("M71").Value = -8681.833500000001
$ws.Range("N71").Value = -23113

$ws.Range("H93").Value = 1797.6154
$ws.Range("I93").Value = 897.3333
$ws.Range("K93").Value = 897.3333
$ws.Range("M93").Value = 350.6667

$ws.Range("H113").Value = 2028.4286
$ws.Range("I113").Value = 2033.1666
$ws.Range("K113").Value = 2033.1666
$ws.Range("M113").Value = 136.8334

$ws = $wb.Worksheets.Item(8)  # WVR
$ws.Range("H107").Value = 1371.5333
$ws.Range("I107").Value = 1316.5454
$ws.Range("J107").Value = 1522.75
$ws.Range("K107").Value = 3949.6362
$ws.Range("L107").Value = 4568.25
$ws.Range("M107").Value = -2029.6362
$ws.Range("N107").Value = -8408.25

$ws.Range("H140").Value = 79899.60000000001
$ws.Range("J140").Value = 82249.5
$ws.Range("L140").Value = 82249.5
$ws.Range("N140").Value = -92609.5

$ws.Range("H141").Value = 93374.75
$ws.Range("J141").Value = 93374.75
$ws.Range("L141").Value = 93374.75
$ws.Range("N141").Value = -103734.75
